$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values (week 2 pbp data)
$ws.Range("B2").Value = 1507.827
$ws.Range("B3").Value = 1464.298
$ws.Range("B7").Value = 1

# Add new formula in F2, matching the style (number format) of C2/D2
$ws.Range("F2").Formula = "=B2-B12"
$ws.Range("F2").NumberFormat = $ws.Range("C2").NumberFormat

# Add new raw data cells
$ws.Range("I2").Value = 1510.923
$ws.Range("J2").Value = 1507.827
$ws.Range("I3").Value = 1509.452
$ws.Range("J3").Value = 1464.298

$wb.Save()
